# Apply citation-numbering edits to section 1 (research motivation) and
# leave the reference-list hyperlinks' visible text unchanged (the diff's
# run-splitting of the hyperlink text is a no-op at the content level —
# the concatenated text of the split runs equals the original text).

$d = $word.ActiveDocument

# --- 1. "...such as education, health, or housing [6]." -> "...[16]." ---
# (Only the in-text citation must change; the "[6] "Alcohol Consumption"," reference
#  list entry a little further down must stay untouched, so we anchor on
#  unique surrounding text.)
$rng = $d.Content
$ok = $rng.Find.Execute(
    "such as education, health, or housing [6].",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "such as education, health, or housing [16].",
    2)
if (-not $ok) { throw "Find/Replace #1 failed" }

# --- 2. "...heavy use of alcohol [4]." -> "...heavy use of alcohol [1]." ---
# (Only the in-text citation; the "[4] "A meta-analysis..." reference list
#  entry must stay untouched.)
$rng = $d.Content
$ok = $rng.Find.Execute(
    "heavy use of alcohol [4].",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "heavy use of alcohol [1].",
    2)
if (-not $ok) { throw "Find/Replace #2 failed" }

# --- 3. "...general population [5]." -> "...general population [1]." ---
$rng = $d.Content
$ok = $rng.Find.Execute(
    "general population [5].",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "general population [1].",
    2)
if (-not $ok) { throw "Find/Replace #3 failed" }

Write-Host "Citation edits applied."
